$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1.1 - Registro de Auditoria e Senhas
# The standalone "SENHA" (password) column is removed from the sample
# table: users now set their own password on first login, and admins
# can only reset (blank) it - there is no stored plaintext password
# to show in this sheet anymore. The ADMIN column shifts from E to D.
# ------------------------------------------------------------------

# --- Header row -----------------------------------------------------
# D1 becomes "ADMIN" (used to be E1); E1 becomes an empty, bordered
# header cell styled like the rest of the header row.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)            # xlPasteFormats
$ws.Range("D1").Value = "ADMIN"

$ws.Range("E1").Interior.ThemeColor = 2        # -> theme 0 (Background 1 / white)
$ws.Range("E1").Borders.Color = 13553360       # RGB(208,206,206) ~ theme 2 tint -0.1
$ws.Range("E1").Borders.LineStyle = 1          # thin
$ws.Range("E1").ClearContents()

# --- Data row 2 -------------------------------------------------------
# D2 becomes the ADMIN boolean flag (used to be E2); E2 becomes an
# empty, bordered cell (no fill) to match the new blank SENHA column.
$ws.Range("E2").Copy()
$ws.Range("D2").PasteSpecial(-4122)            # xlPasteFormats
$ws.Range("D2").Value = $true

$ws.Range("E2").Borders.Color = 13553360       # RGB(208,206,206) ~ theme 2 tint -0.1
$ws.Range("E2").Borders.LineStyle = 1          # thin
$ws.Range("E2").ClearContents()

$ws.Range("B2").Value = "exemplo@exemplo"
$ws.Range("C2").Value = 12345678910

# --- Selection mirrors the author's final cursor position -----------
$ws.Activate()
$ws.Range("B2").Select()
